# Apply the cryptos-list content update (prices / 1h-volume% refresh,
# LEO inserted into the ranking causing rows 23-51 to shift down by one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts plain numeric-looking text (e.g. "1.00", "18.00") to
# real numbers on assignment, which would strip the formatting the sheet
# relies on. Force the affected Price cells to Text first so the literal
# strings round-trip unchanged.
$textCells = @("D5","D6","D8","D13","D18","D19","D21","D23","D24","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

# Row 2
$ws.Range("D2").Value = "62.981.24"
$ws.Range("E2").Value = "  -0.41%  "
# Row 3
$ws.Range("D3").Value = "2.579.74"
$ws.Range("E3").Value = "  +0.42%  "
# Row 4
$ws.Range("E4").Value = "  -0.04%  "
# Row 5
$ws.Range("D5").Value = "584.66"
$ws.Range("E5").Value = "  +0.04%  "
# Row 6
$ws.Range("D6").Value = "144.53"
$ws.Range("E6").Value = "  -2.44%  "
# Row 7
$ws.Range("E7").Value = "  -0.02%  "
# Row 8
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  -1.77%  "
# Row 9
$ws.Range("E9").Value = "  -1.95%  "
# Row 10
$ws.Range("E10").Value = "  -0.90%  "
# Row 11
$ws.Range("E11").Value = "  -0.51%  "
# Row 12
$ws.Range("E12").Value = "  -1.89%  "
# Row 13
$ws.Range("D13").Value = "27.13"
$ws.Range("E13").Value = "  -1.36%  "
# Row 14
$ws.Range("D14").Value = "3.044.24"
$ws.Range("E14").Value = "  +0.53%  "
# Row 15
$ws.Range("D15").Value = "62.880.48"
$ws.Range("E15").Value = "  -0.48%  "
# Row 16
$ws.Range("E16").Value = "  -1.55%  "
# Row 17
$ws.Range("D17").Value = "2.578.56"
$ws.Range("E17").Value = "  +0.38%  "
# Row 18
$ws.Range("D18").Value = "11.11"
$ws.Range("E18").Value = "  -2.18%  "
# Row 19
$ws.Range("D19").Value = "342.69"
$ws.Range("E19").Value = "  -0.01%  "
# Row 20
$ws.Range("E20").Value = "  -1.74%  "
# Row 21
$ws.Range("D21").Value = "6.63"
$ws.Range("E21").Value = "  -3.41%  "
# Row 22
$ws.Range("E22").Value = "  +0.07%  "
# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "5.72"
$ws.Range("E23").Value = "  +1.75%  "
# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "67.69"
$ws.Range("E24").Value = "  +1.50%  "
# Row 25
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("E25").Value = "  +7.30%  "
# Row 26
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "1.60"
$ws.Range("E26").Value = "  -2.40%  "
# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.165"
$ws.Range("E27").Value = "  -3.34%  "
# Row 28
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "7.99"
$ws.Range("E28").Value = "  -2.79%  "
# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.07%  "
# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "8.26"
$ws.Range("E30").Value = "  -2.65%  "
# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.93"
$ws.Range("E31").Value = "  -3.12%  "
# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "460.86"
$ws.Range("E32").Value = "  +0.17%  "
# Row 33
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0801"
$ws.Range("E33").Value = "  -2.89%  "
# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.67"
$ws.Range("E34").Value = "  +2.07%  "
# Row 35
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "176.84"
$ws.Range("E35").Value = "  +0.09%  "
# Row 36
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.07%  "
# Row 37
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "0.400"
$ws.Range("E37").Value = "  -1.66%  "
# Row 38
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "18.86"
$ws.Range("E38").Value = "  -2.02%  "
# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "4.54"
$ws.Range("E39").Value = "  +0.64%  "
# Row 40
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.04%  "
# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  -2.94%  "
# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "158.78"
$ws.Range("E42").Value = "  +4.98%  "
# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "40.01"
$ws.Range("E43").Value = "  +0.78%  "
# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.70"
$ws.Range("E44").Value = "  -3.20%  "
# Row 45
$ws.Range("D45").Value = "21.27"
$ws.Range("E45").Value = "  +0.89%  "
# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.634"
$ws.Range("E46").Value = "  +3.16%  "
# Row 47
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0538"
$ws.Range("E47").Value = "  -2.59%  "
# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.0961"
$ws.Range("E48").Value = "  -2.03%  "
# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0236"
$ws.Range("E49").Value = "  -1.39%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "18.00"
$ws.Range("E50").Value = "  -2.28%  "
# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "11.40"
$ws.Range("E51").Value = "  +0.11%  "
